$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.1748467492742235
$ws0.Range("C2").Value = 0.1120573483709346
$ws0.Range("B3").Value = -0.153701250804759
$ws0.Range("C3").Value = -1.987109540069739
$ws0.Range("B4").Value = -1.420539244477489
$ws0.Range("C4").Value = -0.07219190575164179

$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -0.1517445981469674
$ws1.Range("C2").Value = -0.1559722898806173
$ws1.Range("B3").Value = -1.395836035602749
$ws1.Range("C3").Value = 0.2813611808563983
$ws1.Range("B4").Value = -1.939354137125467
$ws1.Range("C4").Value = 0.2251385441134592
